$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '26.219.96'
Set-TextValue "E2" '  -0.75%  '
Set-TextValue "D3" '1.681.15'
Set-TextValue "D5" '211.35'
Set-TextValue "E5" '  -3.45%  '
Set-TextValue "E6" '  -3.20%  '
Set-TextValue "E7" '  -0.73%  '
Set-TextValue "D8" '0.2683'
Set-TextValue "E8" '  -1.28%  '
Set-TextValue "D9" '0.06309'
Set-TextValue "E9" '  -2.48%  '
Set-TextValue "D10" '21.37'
Set-TextValue "E10" '  -3.18%  '
Set-TextValue "D11" '0.07530'
Set-TextValue "E11" '  -2.35%  '
Set-TextValue "D12" '1.688.42'
Set-TextValue "E12" '  -0.32%  '
Set-TextValue "D13" '4.481'
Set-TextValue "E13" '  -1.38%  '
Set-TextValue "D14" '0.5674'
Set-TextValue "E14" '  -2.55%  '
Set-TextValue "D15" '0.000008141'
Set-TextValue "E15" '  -3.16%  '
Set-TextValue "D16" '66.56'
Set-TextValue "E16" '  +2.02%  '
Set-TextValue "D17" '26.261.66'
Set-TextValue "E17" '  -0.94%  '
Set-TextValue "E18" '  -0.68%  '
Set-TextValue "D19" '4.857'
Set-TextValue "E19" '  -2.02%  '
Set-TextValue "D20" '10.55'
Set-TextValue "E20" '  -4.07%  '
Set-TextValue "D21" '188.78'
Set-TextValue "E21" '  -0.64%  '
Set-TextValue "D22" '6.220'
Set-TextValue "E22" '  -0.23%  '
Set-TextValue "E23" '  -0.69%  '
Set-TextValue "D24" '147.57'
Set-TextValue "E24" '  -1.93%  '
Set-TextValue "D25" '0.1260'
Set-TextValue "E25" '  -3.19%  '
Set-TextValue "D26" '7.629'
Set-TextValue "E26" '  -3.40%  '
Set-TextValue "E27" '  +1.01%  '
Set-TextValue "D28" '0.06450'
Set-TextValue "E28" '  +2.33%  '
Set-TextValue "D29" '1.346'
Set-TextValue "E29" '  -5.34%  '
Set-TextValue "D30" '1.287'
Set-TextValue "E30" '  -3.31%  '
Set-TextValue "E31" '  -1.36%  '
Set-TextValue "D32" '3.490'
Set-TextValue "E32" '  -2.59%  '
Set-TextValue "D33" '1.657'
Set-TextValue "E33" '  -1.16%  '
Set-TextValue "D34" '1.011'
Set-TextValue "E34" '  -2.94%  '
Set-TextValue "D35" '0.6112'
Set-TextValue "E35" '  -1.66%  '
Set-TextValue "E36" '  +0.10%  '
Set-TextValue "D37" '2.718'
Set-TextValue "E37" '  -0.43%  '
Set-TextValue "D38" '6.177'
Set-TextValue "E38" '  -0.68%  '
Set-TextValue "E39" '  -1.28%  '
Set-TextValue "D40" '1.102.18'
Set-TextValue "E40" '  -1.71%  '
Set-TextValue "D41" '0.8689'
Set-TextValue "E41" '  -1.46%  '
Set-TextValue "D43" '100.31'
Set-TextValue "E43" '  -0.89%  '
Set-TextValue "D44" '1.833.10'
Set-TextValue "E44" '  -0.58%  '
Set-TextValue "D45" '0.00000000111'
Set-TextValue "E45" '  +1.59%  '
Set-TextValue "D46" '57.03'
Set-TextValue "E46" '  -0.65%  '
Set-TextValue "D47" '1.002'
Set-TextValue "E47" '  -0.66%  '
Set-TextValue "D48" '0.05266'
Set-TextValue "E48" '  -0.23%  '
Set-TextValue "D49" '7.994'
Set-TextValue "E49" '  -2.53%  '
Set-TextValue "D50" '0.4271'
Set-TextValue "E50" '  -0.84%  '
Set-TextValue "D51" '5.970'
Set-TextValue "E51" '  -1.63%  '
